$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "P_1092"
$ws.Cells.Item(2, 3).Value = 3239.458335206326
$ws.Cells.Item(2, 4).Value = 156.3399505236742
$ws.Cells.Item(2, 5).Value = 91.07593037516521
$ws.Cells.Item(2, 6).Value = 7090.247189282279
$ws.Cells.Item(3, 2).Value = "P_1307"
$ws.Cells.Item(3, 3).Value = 1759.081350000001
$ws.Cells.Item(3, 4).Value = 87.30917999999996
$ws.Cells.Item(3, 5).Value = 53.283533436125
$ws.Cells.Item(3, 6).Value = 3959.599999999998
$ws.Cells.Item(4, 2).Value = "P_1111"
$ws.Cells.Item(4, 3).Value = 1573.940291000001
$ws.Cells.Item(4, 4).Value = 78.49733850000007
$ws.Cells.Item(4, 5).Value = 48.64042827868546
$ws.Cells.Item(4, 6).Value = 3559.970000000003
$ws.Cells.Item(5, 2).Value = "P_1393"
$ws.Cells.Item(5, 3).Value = 2418.835211943
$ws.Cells.Item(5, 4).Value = 122.03978091017
$ws.Cells.Item(5, 5).Value = 52.60679542982393
$ws.Cells.Item(5, 6).Value = 5534.68394150431
$ws.Cells.Item(6, 2).Value = "P_1304"
$ws.Cells.Item(6, 3).Value = 1241.57304
$ws.Cells.Item(6, 4).Value = 61.56360000000002
$ws.Cells.Item(6, 5).Value = 36.57289490219227
$ws.Cells.Item(6, 6).Value = 2792
$ws.Cells.Item(7, 2).Value = "P_1279"
$ws.Cells.Item(7, 3).Value = 3579.051168000003
$ws.Cells.Item(7, 4).Value = 168.40908
$ws.Cells.Item(7, 5).Value = 96.05876624089991
$ws.Cells.Item(7, 6).Value = 7637.599999999999
$ws.Cells.Item(8, 2).Value = "P_1371"
$ws.Cells.Item(8, 3).Value = 1181.371551876699
$ws.Cells.Item(8, 4).Value = 56.64827908777654
$ws.Cells.Item(8, 5).Value = 33.13819390409638
$ws.Cells.Item(8, 6).Value = 2569.0829518266
$ws.Cells.Item(9, 2).Value = "P_1419"
$ws.Cells.Item(9, 3).Value = 1383.360587716136
$ws.Cells.Item(9, 4).Value = 67.53164665920092
$ws.Cells.Item(9, 5).Value = 38.6483330356762
$ws.Cells.Item(9, 6).Value = 3062.659712435416
$ws.Cells.Item(10, 2).Value = "P_1421"
$ws.Cells.Item(10, 3).Value = 1335.709708287825
$ws.Cells.Item(10, 4).Value = 65.61978126510103
$ws.Cells.Item(10, 5).Value = 36.30818143506726
$ws.Cells.Item(10, 6).Value = 2975.953798870795
$ws.Cells.Item(11, 2).Value = "P_1100"
$ws.Cells.Item(11, 3).Value = 3007.23600336858
$ws.Cells.Item(11, 4).Value = 149.7956827499998
$ws.Cells.Item(11, 5).Value = 89.23199304412083
$ws.Cells.Item(11, 6).Value = 6793.454999999992
$ws.Cells.Item(12, 2).Value = "P_1141"
$ws.Cells.Item(12, 3).Value = 3752.842114106415
$ws.Cells.Item(12, 4).Value = 178.5829499999998
$ws.Cells.Item(12, 5).Value = 102.1982032613868
$ws.Cells.Item(12, 6).Value = 8098.999999999988
$ws.Cells.Item(13, 2).Value = "P_1093"
$ws.Cells.Item(13, 3).Value = 2471.975494239886
$ws.Cells.Item(13, 4).Value = 123.3650092500005
$ws.Cells.Item(13, 5).Value = 76.84668256242304
$ws.Cells.Item(13, 6).Value = 5594.785000000023
$ws.Cells.Item(14, 2).Value = "P_1257"
$ws.Cells.Item(14, 3).Value = 1842.183774467457
$ws.Cells.Item(14, 4).Value = 88.96876707959395
$ws.Cells.Item(14, 5).Value = 53.72884433391435
$ws.Cells.Item(14, 6).Value = 4034.864720162991
$ws.Cells.Item(15, 2).Value = "P_1143"
$ws.Cells.Item(15, 3).Value = 2081.658827523183
$ws.Cells.Item(15, 4).Value = 100.4333399999999
$ws.Cells.Item(15, 5).Value = 59.91512060420938
$ws.Cells.Item(15, 6).Value = 4554.799999999996
$ws.Cells.Item(16, 2).Value = "P_1272"
$ws.Cells.Item(16, 3).Value = 3257.068328629367
$ws.Cells.Item(16, 4).Value = 156.9272040000004
$ws.Cells.Item(16, 5).Value = 94.16916893510829
$ws.Cells.Item(16, 6).Value = 7116.880000000015
$ws.Cells.Item(17, 2).Value = "P_1127"
$ws.Cells.Item(17, 3).Value = 4282.020056416666
$ws.Cells.Item(17, 4).Value = 204.9613649999997
$ws.Cells.Item(17, 5).Value = 123.2169825691001
$ws.Cells.Item(17, 6).Value = 9295.299999999983
$ws.Cells.Item(18, 2).Value = "P_1260"
$ws.Cells.Item(18, 3).Value = 3027.787686288529
$ws.Cells.Item(18, 4).Value = 145.7358503319236
$ws.Cells.Item(18, 5).Value = 79.99428898205711
$ws.Cells.Item(18, 6).Value = 6609.335615960253
$ws.Cells.Item(19, 2).Value = "P_1123"
$ws.Cells.Item(19, 3).Value = 3420.729313100619
$ws.Cells.Item(19, 4).Value = 138.0354095644854
$ws.Cells.Item(19, 5).Value = 74.6646107620076
$ws.Cells.Item(19, 6).Value = 6260.109277300924
$ws.Cells.Item(20, 2).Value = "P_1295"
$ws.Cells.Item(20, 3).Value = 1756.763959999997
$ws.Cells.Item(20, 4).Value = 88.51399199999986
$ws.Cells.Item(20, 5).Value = 54.93293609434996
$ws.Cells.Item(20, 6).Value = 4014.239999999993
$ws.Cells.Item(21, 2).Value = "P_1375"
$ws.Cells.Item(21, 3).Value = 2527.974324242419
$ws.Cells.Item(21, 4).Value = 120.1460400000001
$ws.Cells.Item(21, 5).Value = 61.20184795059098
$ws.Cells.Item(21, 6).Value = 5448.800000000005
$ws.Cells.Item(22, 2).Value = "P_1376"
$ws.Cells.Item(22, 3).Value = 1867.287152140784
$ws.Cells.Item(22, 4).Value = 93.68710783880661
$ws.Cells.Item(22, 5).Value = 47.30750649483431
$ws.Cells.Item(22, 6).Value = 4248.848428063791
$ws.Cells.Item(23, 2).Value = "P_1414"
$ws.Cells.Item(23, 3).Value = 918.9812198681035
$ws.Cells.Item(23, 4).Value = 45.86764546122762
$ws.Cells.Item(23, 5).Value = 28.6187456631727
$ws.Cells.Item(23, 6).Value = 2080.165327039801
$ws.Cells.Item(24, 2).Value = "P_1131"
$ws.Cells.Item(24, 3).Value = 2318.926265595423
$ws.Cells.Item(24, 4).Value = 108.7814700000004
$ws.Cells.Item(24, 5).Value = 61.0705177994389
$ws.Cells.Item(24, 6).Value = 4933.400000000017
$ws.Cells.Item(25, 2).Value = "P_1112"
$ws.Cells.Item(25, 3).Value = 2391.252711749999
$ws.Cells.Item(25, 4).Value = 96.61075200000002
$ws.Cells.Item(25, 5).Value = 44.59591799591247
$ws.Cells.Item(25, 6).Value = 4381.440000000001
$ws.Cells.Item(26, 2).Value = "P_1132"
$ws.Cells.Item(26, 3).Value = 1444.018728419883
$ws.Cells.Item(26, 4).Value = 70.53331950000009
$ws.Cells.Item(26, 5).Value = 41.72447031467404
$ws.Cells.Item(26, 6).Value = 3198.790000000004
$ws.Cells.Item(27, 2).Value = "P_1253"
$ws.Cells.Item(27, 3).Value = 1224.50853
$ws.Cells.Item(27, 4).Value = 54.93316500000002
$ws.Cells.Item(27, 5).Value = 35.12202847676254
$ws.Cells.Item(27, 6).Value = 2491.300000000001
$ws.Cells.Item(28, 2).Value = "P_1135"
$ws.Cells.Item(28, 3).Value = 2398.48098589473
$ws.Cells.Item(28, 4).Value = 118.6907400000003
$ws.Cells.Item(28, 5).Value = 71.40030787805004
$ws.Cells.Item(28, 6).Value = 5382.800000000015
$ws.Cells.Item(29, 2).Value = "P_1109"
$ws.Cells.Item(29, 3).Value = 3194.865164999994
$ws.Cells.Item(29, 4).Value = 158.2335562500003
$ws.Cells.Item(29, 5).Value = 98.31126052625729
$ws.Cells.Item(29, 6).Value = 7176.125000000015
$ws.Cells.Item(30, 2).Value = "P_1353"
$ws.Cells.Item(30, 3).Value = 1295.863977020527
$ws.Cells.Item(30, 4).Value = 64.09599683416663
$ws.Cells.Item(30, 5).Value = 40.58505108367701
$ws.Cells.Item(30, 6).Value = 2906.847928987148
$ws.Cells.Item(31, 2).Value = "P_1424"
$ws.Cells.Item(31, 3).Value = 404.1496279225125
$ws.Cells.Item(31, 4).Value = 19.65511432815084
$ws.Cells.Item(31, 5).Value = 12.26175685510235
$ws.Cells.Item(31, 6).Value = 891.3884049048
$ws.Cells.Item(32, 2).Value = "P_1427"
$ws.Cells.Item(32, 3).Value = 1993.503319245416
$ws.Cells.Item(32, 4).Value = 98.43651703802752
$ws.Cells.Item(32, 5).Value = 58.59265312081038
$ws.Cells.Item(32, 6).Value = 4464.241135511452
$ws.Cells.Item(33, 2).Value = "P_1217"
$ws.Cells.Item(33, 3).Value = 614.8600797368414
$ws.Cells.Item(33, 4).Value = 25.474365
$ws.Cells.Item(33, 5).Value = 11.32326814729268
$ws.Cells.Item(33, 6).Value = 1155.3
$ws.Cells.Item(34, 2).Value = "P_1271"
$ws.Cells.Item(34, 3).Value = 4166.875983696955
$ws.Cells.Item(34, 4).Value = 190.32678
$ws.Cells.Item(34, 5).Value = 100.2043873143121
$ws.Cells.Item(34, 6).Value = 8631.599999999999
$ws.Cells.Item(35, 2).Value = "P_1368"
$ws.Cells.Item(35, 3).Value = 1519.11522584357
$ws.Cells.Item(35, 4).Value = 74.43648098050134
$ws.Cells.Item(35, 5).Value = 46.05181330529004
$ws.Cells.Item(35, 6).Value = 3375.804126099834
$ws.Cells.Item(36, 2).Value = "P_1094"
$ws.Cells.Item(36, 3).Value = 2603.176699124993
$ws.Cells.Item(36, 4).Value = 114.314256
$ws.Cells.Item(36, 5).Value = 66.33131540430304
$ws.Cells.Item(36, 6).Value = 5184.319999999997
$ws.Cells.Item(37, 2).Value = "P_1285"
$ws.Cells.Item(37, 3).Value = 1783.073519999994
$ws.Cells.Item(37, 4).Value = 86.91228000000004
$ws.Cells.Item(37, 5).Value = 50.52206982505907
$ws.Cells.Item(37, 6).Value = 3941.600000000001
$ws.Cells.Item(38, 2).Value = "P_1288"
$ws.Cells.Item(38, 3).Value = 783.5702944809786
$ws.Cells.Item(38, 4).Value = 38.28530474999997
$ws.Cells.Item(38, 5).Value = 21.60099836996058
$ws.Cells.Item(38, 6).Value = 1736.294999999999
$ws.Cells.Item(39, 2).Value = "P_1357"
$ws.Cells.Item(39, 3).Value = 2013.695720089836
$ws.Cells.Item(39, 4).Value = 102.2198741718037
$ws.Cells.Item(39, 5).Value = 65.1908425176393
$ws.Cells.Item(39, 6).Value = 4635.821957904929
$ws.Cells.Item(40, 2).Value = "P_1280"
$ws.Cells.Item(40, 3).Value = 1363.089244363638
$ws.Cells.Item(40, 4).Value = 66.38461199999999
$ws.Cells.Item(40, 5).Value = 39.61473453095603
$ws.Cells.Item(40, 6).Value = 3010.639999999999
$ws.Cells.Item(41, 2).Value = "P_1281"
$ws.Cells.Item(41, 3).Value = 2411.072585192307
$ws.Cells.Item(41, 4).Value = 117.436095
$ws.Cells.Item(41, 5).Value = 68.50094724358814
$ws.Cells.Item(41, 6).Value = 5325.900000000001
$ws.Cells.Item(42, 2).Value = "P_1296"
$ws.Cells.Item(42, 3).Value = 1841.298604495516
$ws.Cells.Item(42, 4).Value = 88.03021499999994
$ws.Cells.Item(42, 5).Value = 52.77573966803654
$ws.Cells.Item(42, 6).Value = 3992.299999999997
$ws.Cells.Item(43, 2).Value = "P_1367"
$ws.Cells.Item(43, 3).Value = 1683.73791217941
$ws.Cells.Item(43, 4).Value = 83.50570956409045
$ws.Cells.Item(43, 5).Value = 50.09945291411204
$ws.Cells.Item(43, 6).Value = 3787.107009709317
$ws.Cells.Item(44, 2).Value = "P_1125a"
$ws.Cells.Item(44, 3).Value = 1017.167607194544
$ws.Cells.Item(44, 4).Value = 41.83201878493561
$ws.Cells.Item(44, 5).Value = 20.95130077938576
$ws.Cells.Item(44, 6).Value = 1897.143709067375
$ws.Cells.Item(45, 2).Value = "P_1087"
$ws.Cells.Item(45, 3).Value = 3704.786575062497
$ws.Cells.Item(45, 4).Value = 181.5042442499993
$ws.Cells.Item(45, 5).Value = 108.871073795229
$ws.Cells.Item(45, 6).Value = 8231.484999999966
$ws.Cells.Item(46, 2).Value = "P_1422"
$ws.Cells.Item(46, 3).Value = 992.9109013134209
$ws.Cells.Item(46, 4).Value = 49.29635530629912
$ws.Cells.Item(46, 5).Value = 28.5739729539983
$ws.Cells.Item(46, 6).Value = 2235.66237216776
$ws.Cells.Item(47, 2).Value = "P_1137"
$ws.Cells.Item(47, 3).Value = 1736.542789545504
$ws.Cells.Item(47, 4).Value = 78.35881930843867
$ws.Cells.Item(47, 5).Value = 41.2342618612266
$ws.Cells.Item(47, 6).Value = 3553.687950496084
$ws.Cells.Item(48, 2).Value = "P_1294"
$ws.Cells.Item(48, 3).Value = 788.6366614535214
$ws.Cells.Item(48, 4).Value = 35.58278763351645
$ws.Cells.Item(48, 5).Value = 18.68265108277637
$ws.Cells.Item(48, 6).Value = 1613.731865465599
$ws.Cells.Item(49, 2).Value = "P_1091"
$ws.Cells.Item(49, 3).Value = 6543.751701229684
$ws.Cells.Item(49, 4).Value = 325.5459794999994
$ws.Cells.Item(49, 5).Value = 201.0975713303817
$ws.Cells.Item(49, 6).Value = 14763.98999999997
$ws.Cells.Item(50, 2).Value = "P_1125"
$ws.Cells.Item(50, 3).Value = 910.5217470224015
$ws.Cells.Item(50, 4).Value = 37.52948423153453
$ws.Cells.Item(50, 5).Value = 21.83733499788648
$ws.Cells.Item(50, 6).Value = 1702.017425466418
$ws.Cells.Item(51, 2).Value = "P_1098"
$ws.Cells.Item(51, 3).Value = 2622.280799520234
$ws.Cells.Item(51, 4).Value = 121.5743120129884
$ws.Cells.Item(51, 5).Value = 58.89065695678865
$ws.Cells.Item(51, 6).Value = 5513.574240951853
